$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.013591766357422
$ws.Range("B1").Value = 3.701288461685181
$ws.Range("C1").Value = 3.911319017410278
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 3.436680316925049
